$wb = $excel.ActiveWorkbook
$heating = $wb.Worksheets.Item("HEATING")
$cooling = $wb.Worksheets.Item("COOLING")

# Move the "district heating - natural gas-fired boiler" row from HEATING
# to COOLING (as a "district cooling" variant used for the absorption
# chiller), fixing the LCA database linkage.
$heating.Range("A5:I5").Copy($cooling.Range("A7:I7")) | Out-Null

# Update the moved row's text to reflect the cooling-side entry.
$cooling.Range("A7").Value = "district cooling - natural gas-fired boiler for absorption chiller"
$cooling.Range("C7").Value = "NG"

# Remove the now-duplicated row from HEATING.
$heating.Rows.Item(5).Delete() | Out-Null

# Restore the selection on HEATING left over from the row deletion.
$heating.Range("A5:XFD5").Select() | Out-Null

# COOLING becomes the active sheet, with the selection left on the newly
# added row.
$cooling.Activate() | Out-Null
$cooling.Range("D26").Select() | Out-Null
